# Generate Report for Handoff
# Updates the status/date for the 96ea4119-...-md file now that it is
# "Ready for handoff" (instead of "Handed back: in sync with en-US"),
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# Row 3 corresponds to file 96ea4119-27f9-41b6-9d33-b6b848f64680.md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(3, 2).Value = "Ready for handoff"          # B3: zh-cn status
$wsOverview.Cells.Item(3, 3).Value = "Ready for handoff"          # C3: de-de status
$wsOverview.Cells.Item(3, 4).Value = "2016-03-21 12:40:54"        # D3: Latest Handoff Date

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Cells.Item(3, 3).Value = "Ready for handoff"               # C3: Status
$wsZhCn.Cells.Item(3, 5).Value = "2016-03-21 12:40:51"             # E3: Latest Handoff Datetime

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Cells.Item(3, 3).Value = "Ready for handoff"                # C3: Status
$wsDeDe.Cells.Item(3, 5).Value = "2016-03-21 12:40:54"              # E3: Latest Handoff Datetime
